$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row to the table "Condicion_Pacientes" (data for 2020-05-17)
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.ListRows.Add() | Out-Null

# Copy the formatting from the previous data row (row 65) onto the new row (66)
$ws.Range("A65:F65").Copy()
$ws.Range("A66:F66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values
$ws.Range("A66").Value = 43968
$ws.Range("B66").Value = 366
$ws.Range("C66").Value = 81
$ws.Range("D66").Value = 252
$ws.Range("E66").Value = 14
$ws.Range("F66").Value = 17

# Update the active cell / selection to reflect the last edited cell
$ws.Range("C66").Select()
